$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("F1").Value = "status"
$ws.Range("G1").Value = "statusId"
$ws.Range("H1").Value = "country"
$ws.Range("I1").Value = "city"

# Row 4
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "gger"
$ws.Range("C4").Value = "Realsdasd"
$ws.Range("D4").Value = "asda"
$ws.Range("E4").Value = "2021-09-06T21:32:47.954Z"
$ws.Range("F4").Value = "Sent"
$ws.Range("G4").Value = 0

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "gger"
$ws.Range("C5").Value = "Realsdasd"
$ws.Range("D5").Value = "asda"
$ws.Range("E5").Value = "2021-09-06T21:35:48.008Z"
$ws.Range("F5").Value = "Sent"
$ws.Range("G5").Value = 0

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "gger"
$ws.Range("C6").Value = "Realsdasd"
$ws.Range("D6").Value = "asda"
$ws.Range("E6").Value = "2021-09-06T21:36:51.238Z"
$ws.Range("F6").Value = "Sent"
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = "canada"
$ws.Range("I6").Value = "toronto"

# Keep the "numbers stored as text" ignored-error marker covering the
# full used range (it originally covered A1:E3, matching the dimension).
$r = $ws.Range("A1:I6")
$r.Errors.Item([Microsoft.Office.Interop.Excel.XlErrorChecks]::xlNumberAsText).Ignore = $true
